$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; this pushes the existing rows 48:68
# down to 49:69 (dimension grows from T68 to T69).
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Range("A48").Value = 6
$ws.Range("B48").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44574
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100101
$ws.Range("H48").Value = "Berries"
$ws.Range("I48").Value = 100101008
$ws.Range("J48").Value = "Mora"
$ws.Range("K48").Value = "Sin especificar"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 350
$ws.Range("N48").Value = 6000
$ws.Range("O48").Value = 6000
$ws.Range("P48").Value = 6000
$ws.Range("Q48").Value = "$/bandeja 2 kilos"
$ws.Range("R48").Value = "Provincia de Linares"
$ws.Range("S48").Value = 3000
$ws.Range("T48").Value = 2
